$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4059.121
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 4059.121
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 12177.363
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -12513.363

$ws.Range("H41").Value = 1114.8
$ws.Range("J41").Value = 631.6
$ws.Range("L41").Value = 631.6
$ws.Range("N41").Value = -1511.6

$ws.Range("H52").Value = 35000
$ws.Range("I52").Value = 27500
$ws.Range("K52").Value = 82500
$ws.Range("M52").Value = -82340

$ws.Range("H62").Value = 10337.556
$ws.Range("I62").Value = 9334
$ws.Range("J62").Value = 12344.667
$ws.Range("K62").Value = 9334
$ws.Range("L62").Value = 12344.667
$ws.Range("M62").Value = -8710
$ws.Range("N62").Value = -13592.667

$ws.Range("H65").Value = 10337.556
$ws.Range("I65").Value = 9334
$ws.Range("J65").Value = 12344.667
$ws.Range("K65").Value = 46670
$ws.Range("L65").Value = 61723.335
$ws.Range("M65").Value = -43550
$ws.Range("N65").Value = -67963.33499999999

$ws.Range("H86").Value = 250085140
$ws.Range("J86").Value = 166835980
$ws.Range("L86").Value = 166835980
$ws.Range("N86").Value = -166838226

$ws.Range("H89").Value = 250085140
$ws.Range("J89").Value = 166835980
$ws.Range("L89").Value = 834179900
$ws.Range("N89").Value = -834191132

$ws.Range("H116").Value = 19036.1
$ws.Range("I116").Value = 19036.1
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 19036.1
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -15594.1
$ws.Range("N116").ClearContents()

$ws.Range("H132").Value = 2309.6
$ws.Range("I132").Value = 2586.0588
$ws.Range("K132").Value = 7758.176399999999
$ws.Range("M132").Value = -5228.176399999999

$ws.Range("H137").Value = 2369.7827
$ws.Range("I137").Value = 2969.6924
$ws.Range("J137").Value = 1589.9
$ws.Range("K137").Value = 8909.0772
$ws.Range("L137").Value = 4769.700000000001
$ws.Range("M137").Value = -6359.0772
$ws.Range("N137").Value = -9869.700000000001

$ws.Range("H138").Value = 3119.3635
$ws.Range("I138").Value = 1394.2122
$ws.Range("J138").Value = 3981.9395
$ws.Range("K138").Value = 4182.6366
$ws.Range("L138").Value = 11945.8185
$ws.Range("M138").Value = 957.3634000000002
$ws.Range("N138").Value = -22225.8185

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4498.8
$ws.Range("I45").Value = 3910.353
$ws.Range("K45").Value = 3910.353
$ws.Range("M45").Value = -3533.353

$ws.Range("H61").Value = 4229.923
$ws.Range("I61").Value = 4300
$ws.Range("J61").Value = 4117.8
$ws.Range("K61").Value = 4300
$ws.Range("L61").Value = 4117.8
$ws.Range("M61").Value = -4088
$ws.Range("N61").Value = -4541.8

$ws.Range("H63").Value = 4573.636
$ws.Range("J63").Value = 7400
$ws.Range("L63").Value = 7400
$ws.Range("N63").Value = -8772

$ws.Range("H66").Value = 4573.636
$ws.Range("J66").Value = 7400
$ws.Range("L66").Value = 37000
$ws.Range("N66").Value = -43864

$ws.Range("H102").Value = 1421.5807
$ws.Range("I102").Value = 1198.84
$ws.Range("J102").Value = 2349.6667
$ws.Range("K102").Value = 1198.84
$ws.Range("L102").Value = 2349.6667
$ws.Range("M102").Value = 423.1600000000001
$ws.Range("N102").Value = -5593.6667

$ws.Range("H110").Value = 1923.6111
$ws.Range("I110").Value = 738.7273
$ws.Range("J110").Value = 3785.5715
$ws.Range("K110").Value = 738.7273
$ws.Range("L110").Value = 3785.5715
$ws.Range("M110").Value = 1306.2727
$ws.Range("N110").Value = -7875.5715

$ws.Range("H132").Value = 1659.8474
$ws.Range("I132").Value = 1526.1091
$ws.Range("K132").Value = 4578.3273
$ws.Range("M132").Value = -2048.3273

$ws.Range("H136").Value = 4229.923
$ws.Range("I136").Value = 4300
$ws.Range("J136").Value = 4117.8
$ws.Range("K136").Value = 12900
$ws.Range("L136").Value = 12353.4
$ws.Range("M136").Value = -10350
$ws.Range("N136").Value = -17453.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 185.84616
$ws.Range("I22").Value = 170.18182
$ws.Range("J22").Value = 272
$ws.Range("K22").Value = 170.18182
$ws.Range("L22").Value = 272
$ws.Range("M22").Value = 2.818180000000012
$ws.Range("N22").Value = -618

$ws.Range("H105").Value = 2655.4375
$ws.Range("I105").Value = 2576
$ws.Range("K105").Value = 2576
$ws.Range("M105").Value = -829

$ws.Range("H134").Value = 5132978
$ws.Range("I134").Value = 6670421.5
$ws.Range("K134").Value = 20011264.5
$ws.Range("M134").Value = -20008729.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1645.3966
$ws.Range("I31").Value = 916.375
$ws.Range("J31").Value = 2160
$ws.Range("K31").Value = 916.375
$ws.Range("L31").Value = 2160
$ws.Range("M31").Value = -621.375
$ws.Range("N31").Value = -2750

$ws.Range("H34").Value = 1645.3966
$ws.Range("I34").Value = 916.375
$ws.Range("J34").Value = 2160
$ws.Range("K34").Value = 916.375
$ws.Range("L34").Value = 2160
$ws.Range("M34").Value = -714.375
$ws.Range("N34").Value = -2564

$ws.Range("H43").Value = 35450
$ws.Range("J43").Value = 35450
$ws.Range("L43").Value = 35450
$ws.Range("N43").Value = -35818

$ws.Range("H58").Value = 3278.5925
$ws.Range("I58").Value = 1892.5
$ws.Range("J58").Value = 4093.9412
$ws.Range("K58").Value = 1892.5
$ws.Range("L58").Value = 4093.9412
$ws.Range("M58").Value = -1689.5
$ws.Range("N58").Value = -4499.9412

$ws.Range("H99").Value = 2973.1875
$ws.Range("I99").Value = 2758.6
$ws.Range("J99").Value = 3330.8333
$ws.Range("K99").Value = 2758.6
$ws.Range("L99").Value = 3330.8333
$ws.Range("M99").Value = -1260.6
$ws.Range("N99").Value = -6326.8333

$ws.Range("H101").Value = 35450
$ws.Range("J101").Value = 35450
$ws.Range("L101").Value = 35450
$ws.Range("N101").Value = -41940

$ws.Range("H126").Value = 2973.1875
$ws.Range("I126").Value = 2758.6
$ws.Range("J126").Value = 3330.8333
$ws.Range("K126").Value = 8275.799999999999
$ws.Range("L126").Value = 9992.499899999999
$ws.Range("M126").Value = -5805.799999999999
$ws.Range("N126").Value = -14932.4999

$ws.Range("H132").Value = 4801.125
$ws.Range("I132").Value = 4646.909
$ws.Range("K132").Value = 13940.727
$ws.Range("M132").Value = -11410.727

$ws.Range("H134").Value = 5628.5
$ws.Range("I134").Value = 5628.5
$ws.Range("K134").Value = 16885.5
$ws.Range("M134").Value = -14350.5

$ws.Range("H136").Value = 3278.5925
$ws.Range("I136").Value = 1892.5
$ws.Range("J136").Value = 4093.9412
$ws.Range("K136").Value = 5677.5
$ws.Range("L136").Value = 12281.8236
$ws.Range("M136").Value = -3127.5
$ws.Range("N136").Value = -17381.8236

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1236.0454
$ws.Range("J107").Value = 1504.5
$ws.Range("L107").Value = 4513.5
$ws.Range("N107").Value = -8353.5

$ws.Range("H129").Value = 2428.5908
$ws.Range("J129").Value = 2499.9524
$ws.Range("L129").Value = 7499.8572
$ws.Range("N129").Value = -17499.8572

$ws.Range("H131").Value = 2412.1428
$ws.Range("J131").Value = 1984.826
$ws.Range("L131").Value = 5954.478
$ws.Range("N131").Value = -16034.478

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2881.3333
$ws.Range("I102").Value = 2556.5386
$ws.Range("J102").Value = 4992.5
$ws.Range("K102").Value = 2556.5386
$ws.Range("L102").Value = 4992.5
$ws.Range("M102").Value = -934.5385999999999
$ws.Range("N102").Value = -8236.5

$ws.Range("H103").Value = 84214.86
$ws.Range("J103").Value = 84214.86
$ws.Range("L103").Value = 84214.86
$ws.Range("N103").Value = -86558.86

$ws.Range("H132").Value = 2467.0833
$ws.Range("I132").Value = 1572.1428
$ws.Range("J132").Value = 3720
$ws.Range("K132").Value = 4716.428400000001
$ws.Range("L132").Value = 11160
$ws.Range("M132").Value = -2186.428400000001
$ws.Range("N132").Value = -16220

$ws.Range("H134").Value = 54975
$ws.Range("J134").Value = 54975
$ws.Range("L134").Value = 164925
$ws.Range("N134").Value = -169995

$ws.Range("H135").Value = 227698.33
$ws.Range("J135").Value = 227698.33
$ws.Range("L135").Value = 227698.33
$ws.Range("N135").Value = -237838.33

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3808.3333
$ws.Range("I7").Value = 3792.5
$ws.Range("K7").Value = 3792.5
$ws.Range("M7").Value = -3680.5

$ws.Range("H22").Value = 1723.3077
$ws.Range("I22").Value = 1243.55
$ws.Range("K22").Value = 1243.55
$ws.Range("M22").Value = -948.55

$ws.Range("H27").Value = 1723.3077
$ws.Range("I27").Value = 1243.55
$ws.Range("K27").Value = 1243.55
$ws.Range("M27").Value = -1136.55

$ws.Range("H40").Value = 2866.6667
$ws.Range("I40").Value = 2453.7273
$ws.Range("K40").Value = 2453.7273
$ws.Range("M40").Value = -2317.7273

$ws.Range("H126").Value = 3808.3333
$ws.Range("I126").Value = 3792.5
$ws.Range("K126").Value = 11377.5
$ws.Range("M126").Value = -8907.5

$ws.Range("H132").Value = 12556.765
$ws.Range("I132").Value = 11565
$ws.Range("K132").Value = 34695
$ws.Range("M132").Value = -32165

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 29408.5
$ws.Range("J45").Value = 29408.5
$ws.Range("L45").Value = 29408.5
$ws.Range("N45").Value = -30390.5

$ws.Range("H132").Value = 2198.2917
$ws.Range("I132").Value = 2250.5217
$ws.Range("K132").Value = 6751.5651
$ws.Range("M132").Value = -4221.5651

$ws.Range("H136").Value = 49716.453
$ws.Range("I136").Value = 3854.818
$ws.Range("K136").Value = 11564.454
$ws.Range("M136").Value = -9014.454000000002
